# Regenerate merged AHB files
# - Rename the "_old" header columns to "_FV2410" and the "_new" header
#   columns to "_FV2504" (A1:J1 and L1:U1, K1 stays "diff").
# - Freeze the header row (row 1).
# - Turn the A1:U58 range into an Excel Table ("Table1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fields = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $fields.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $fields[$i] + "_FV2410"
    $ws.Range($newCols[$i] + "1").Value = $fields[$i] + "_FV2504"
}

# Freeze the header row.
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)

# Convert the data range into a table.
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
[void]($tbl.TableStyle = "")
